$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:A6").NumberFormat = "@"

$ws.Range("A1").Value = "542324853269856256"
$ws.Range("A2").Value = "458987644651700224"
$ws.Range("A3").Value = "295323105994342401"
$ws.Range("A4").Value = "617313325671383051"
$ws.Range("A5").Value = "721016295646298142"
$ws.Range("A6").Value = "649404929462501386"
